$d = $word.ActiveDocument

# 1. Update the table caption / footnote text describing the proportion metric.
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$null = $find.Execute(
    " Proportion of total biomass as below ground (roots) or above ground (shoots) biomass.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    " Proportion of total plant biomass as above ground (shoot) or above ground (root) biomass.",
    2)

# 2. Append units " (g/g)" to the "Shoot biomass" row label.
$rngShoot = $d.Content
$findShoot = $rngShoot.Find
$findShoot.ClearFormatting()
$findShoot.Text = "Shoot biomass"
$foundShoot = $findShoot.Execute()
if ($foundShoot) {
    $rngShoot.Collapse(0)
    $rngShoot.InsertAfter(" (g/g)")
}

# 3. Append units " (g/g)" to the "Root biomass" row label.
$rngRoot = $d.Content
$findRoot = $rngRoot.Find
$findRoot.ClearFormatting()
$findRoot.Text = "Root biomass"
$foundRoot = $findRoot.Execute()
if ($foundRoot) {
    $rngRoot.Collapse(0)
    $rngRoot.InsertAfter(" (g/g)")
}
